$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234; all existing rows from 234 downward
# shift down by one (234 -> 235, ..., 274 -> 275).
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new weekly price record.
$ws.Cells.Item(234, 1).Value = 6
$ws.Cells.Item(234, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(234, 3).Value = "Metropolitana"
$ws.Cells.Item(234, 4).Value = 44491
$ws.Cells.Item(234, 5).Value = 13
$ws.Cells.Item(234, 6).Value = 100112030
$ws.Cells.Item(234, 7).Value = "Poroto granado"
$ws.Cells.Item(234, 8).Value = "Sin especificar"
$ws.Cells.Item(234, 9).Value = "Primera"
$ws.Cells.Item(234, 10).Value = 380
$ws.Cells.Item(234, 11).Value = 30000
$ws.Cells.Item(234, 12).Value = 35000
$ws.Cells.Item(234, 13).Value = 31974
$ws.Cells.Item(234, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(234, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(234, 16).Value = 1279
$ws.Cells.Item(234, 17).Value = 25
$ws.Cells.Item(234, 18).Value = "Hortaliza"

# Make sure the format of the date cell matches the rest of column D
# (date-formatted cells use a dedicated date number format).
$ws.Cells.Item(234, 4).NumberFormat = $ws.Cells.Item(235, 4).NumberFormat()
